$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 3.483499999999999
$ws.Range("C2").Value = 1.0008
$ws.Range("D2").Value = 0.09755
$ws.Range("E2").Value = 27.37635
$ws.Range("F2").Value = 25.67095
$ws.Range("G2").Value = 25.67095
$ws.Range("K2").Value = 63.2648
$ws.Range("L2").Value = 25.671
$ws.Range("M2").Value = 37.5938
$ws.Range("N2").Value = 3.4142
$ws.Range("O2").Value = 34.17960000000001
$ws.Range("B3").Value = 7.121
$ws.Range("C3").Value = 4.587
$ws.Range("E3").Value = 36.987
$ws.Range("F3").Value = 27.605
$ws.Range("G3").Value = 27.4
$ws.Range("H3").Value = 0.205
$ws.Range("I3").Value = 0.205
$ws.Range("K3").Value = 28.7166
$ws.Range("L3").Value = 27.4
$ws.Range("M3").Value = 1.3166
$ws.Range("N3").Value = 1.3166
$ws.Range("B4").Value = 7.477
$ws.Range("C4").Value = 4.564
$ws.Range("E4").Value = 45.696
$ws.Range("F4").Value = 27.69
$ws.Range("G4").Value = 27.68970684931507
$ws.Range("K4").Value = 29.0054
$ws.Range("L4").Value = 27.69
$ws.Range("M4").Value = 1.3154
$ws.Range("N4").Value = 1.3154

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 3.29145
$ws.Range("C2").Value = 1.10835
$ws.Range("D2").Value = 1.74335
$ws.Range("E2").Value = 24.98080000000001
$ws.Range("F2").Value = 25.60595
$ws.Range("G2").Value = 25.60595
$ws.Range("K2").Value = 74.1794
$ws.Range("L2").Value = 25.606
$ws.Range("M2").Value = 48.57340000000001
$ws.Range("N2").Value = 6.6616
$ws.Range("O2").Value = 41.9118
$ws.Range("B3").Value = 7.11
$ws.Range("C3").Value = 5.806
$ws.Range("D3").Value = 10.287
$ws.Range("E3").Value = 32.819
$ws.Range("F3").Value = 27.856
$ws.Range("G3").Value = 27.651
$ws.Range("H3").Value = 0.204
$ws.Range("I3").Value = 0.204
$ws.Range("K3").Value = 27.651
$ws.Range("L3").Value = 27.651
$ws.Range("B4").Value = 7.426
$ws.Range("C4").Value = 10.786
$ws.Range("E4").Value = 24.472
$ws.Range("F4").Value = 27.895
$ws.Range("G4").Value = 27.89524931506849
$ws.Range("K4").Value = 27.895
$ws.Range("L4").Value = 27.895

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 3.24255
$ws.Range("C2").Value = 0.9843
$ws.Range("D2").Value = 1.66455
$ws.Range("E2").Value = 24.4588
$ws.Range("F2").Value = 25.56495
$ws.Range("G2").Value = 25.56495
$ws.Range("K2").Value = 128.2812
$ws.Range("L2").Value = 25.565
$ws.Range("M2").Value = 102.7162
$ws.Range("N2").Value = 4.5192
$ws.Range("O2").Value = 98.197
$ws.Range("B3").Value = 7.11
$ws.Range("C3").Value = 5.806
$ws.Range("D3").Value = 10.287
$ws.Range("E3").Value = 32.819
$ws.Range("F3").Value = 27.856
$ws.Range("G3").Value = 27.651
$ws.Range("H3").Value = 0.204
$ws.Range("I3").Value = 0.204
$ws.Range("K3").Value = 27.651
$ws.Range("L3").Value = 27.651
$ws.Range("B4").Value = 7.426
$ws.Range("C4").Value = 10.786
$ws.Range("E4").Value = 24.472
$ws.Range("F4").Value = 27.895
$ws.Range("G4").Value = 27.89524931506849
$ws.Range("K4").Value = 27.895
$ws.Range("L4").Value = 27.895

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 3.245299999999999
$ws.Range("C2").Value = 1.10935
$ws.Range("D2").Value = 1.7459
$ws.Range("E2").Value = 23.16885
$ws.Range("F2").Value = 25.5542
$ws.Range("G2").Value = 25.5542
$ws.Range("K2").Value = 54.93060000000001
$ws.Range("L2").Value = 25.554
$ws.Range("M2").Value = 29.3766
$ws.Range("N2").Value = 6.3892
$ws.Range("O2").Value = 22.9876
$ws.Range("B3").Value = 7.052
$ws.Range("C3").Value = 4.591
$ws.Range("D3").Value = 12.8
$ws.Range("E3").Value = 35.315
$ws.Range("F3").Value = 27.852
$ws.Range("G3").Value = 27.608
$ws.Range("H3").Value = 0.244
$ws.Range("I3").Value = 0.244
$ws.Range("K3").Value = 27.7508
$ws.Range("L3").Value = 27.607
$ws.Range("M3").Value = 0.1438
$ws.Range("N3").Value = 0.1438
$ws.Range("B4").Value = 7.426
$ws.Range("C4").Value = 10.786
$ws.Range("E4").Value = 24.472
$ws.Range("F4").Value = 27.895
$ws.Range("G4").Value = 27.89524931506849
$ws.Range("K4").Value = 27.895
$ws.Range("L4").Value = 27.895

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 3.1206
$ws.Range("C2").Value = 0.93
$ws.Range("D2").Value = 1.73225
$ws.Range("E2").Value = 21.8446
$ws.Range("F2").Value = 25.46815
$ws.Range("G2").Value = 25.46815
$ws.Range("K2").Value = 174.9578
$ws.Range("L2").Value = 25.468
$ws.Range("M2").Value = 149.4898
$ws.Range("N2").Value = 13.7238
$ws.Range("O2").Value = 135.766
$ws.Range("B3").Value = 5.839
$ws.Range("C3").Value = 11.021
$ws.Range("E3").Value = 15.362
$ws.Range("F3").Value = 27.229
$ws.Range("G3").Value = 27.229
$ws.Range("K3").Value = 31.4578
$ws.Range("L3").Value = 27.229
$ws.Range("M3").Value = 4.228800000000001
$ws.Range("N3").Value = 4.228800000000001
$ws.Range("B4").Value = 5.839
$ws.Range("C4").Value = 11.021
$ws.Range("E4").Value = 15.362
$ws.Range("F4").Value = 27.229
$ws.Range("G4").Value = 27.22919178082192
$ws.Range("K4").Value = 31.4578
$ws.Range("L4").Value = 27.229
$ws.Range("M4").Value = 4.228800000000001
$ws.Range("N4").Value = 4.228800000000001
